$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.224.04'
$ws.Range("E2").Value = '  +0.29%  '

$ws.Range("D3").Value = '2.616.79'
$ws.Range("E3").Value = '  +0.01%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '520.61'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.84%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.37'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.02%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.19%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.568'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.47%  '

$ws.Range("D9").Value = '2.615.54'
$ws.Range("E9").Value = '  -0.19%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.68'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.23%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.102'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.59%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.327'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.48%  '

$ws.Range("E13").Value = '  -0.87%  '

$ws.Range("D14").Value = '3.070.35'
$ws.Range("E14").Value = '  -0.29%  '

$ws.Range("D15").Value = '58.192.80'
$ws.Range("E15").Value = '  +0.16%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.50'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.55%  '

$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000134'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.26%  '

$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '2.609.12'
$ws.Range("E18").Value = '  -0.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '340.51'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.94%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.34'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.89%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.31'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.25%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.38'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.14%  '

$ws.Range("E23").Value = '  -0.17%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.28'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.32%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.166'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.09%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.403'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.58%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.996'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.39%  '

$ws.Range("D28").Value = '2.707.64'
$ws.Range("E28").Value = '  -1.46%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.04'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.20%  '

$ws.Range("D30").Value = '0.0₃0754'
$ws.Range("E30").Value = '  -3.30%  '

$ws.Range("E31").Value = '  -0.07%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.27'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.98%  '

$ws.Range("E33").Value = '  +0.57%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.80'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.73%  '

$ws.Range("E35").Value = '  -1.55%  '

$ws.Range("E36").Value = '  -1.00%  '

$ws.Range("E37").Value = '  -1.50%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.879'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.02%  '

$ws.Range("E39").Value = '  +2.96%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.844'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.16%  '

$ws.Range("E41").Value = '  -1.44%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.57'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.77%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.998'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.23%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '274.51'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.52%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.598'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.02%  '

$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0956'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.84%  '

$ws.Range("B47").Value = 'WhiteBITCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.66'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.44%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.90'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.40%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0523'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.50%  '

$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.05'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.85%  '

$ws.Range("D51").Value = '1.982.52'
$ws.Range("E51").Value = '  -2.26%  '
